$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "34.590.43"
$ws.Range("E2").Value = "  +1.45%  "
# Row 3
$ws.Range("D3").Value = "1.830.51"
$ws.Range("E3").Value = "  +2.62%  "
# Row 4
$ws.Range("E4").Value = "  -0.50%  "
# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "225.01"
$ws.Range("E5").Value = "  -0.60%  "
# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.554"
$ws.Range("E6").Value = "  +0.69%  "
# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.996"
$ws.Range("E7").Value = "  -0.59%  "
# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "32.91"
$ws.Range("E8").Value = "  +4.86%  "
# Row 9
$ws.Range("E9").Value = "  +3.71%  "
# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0707"
$ws.Range("E10").Value = "  +7.41%  "
# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0930"
$ws.Range("E11").Value = "  -0.13%  "
# Row 12
$ws.Range("D12").Value = "2.077.30"
$ws.Range("E12").Value = "  +1.61%  "
# Row 13
$ws.Range("B13").Value = "Chainlink"
$ws.Range("C13").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "11.19"
$ws.Range("E13").Value = "  -0.36%  "
# Row 14
$ws.Range("B14").Value = "WrappedEther"
$ws.Range("C14").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D14").Value = "1.813.39"
$ws.Range("E14").Value = "  +1.15%  "
# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.652"
$ws.Range("E15").Value = "  +3.85%  "
# Row 16
$ws.Range("D16").Value = "34.608.61"
$ws.Range("E16").Value = "  +1.48%  "
# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "4.33"
$ws.Range("E17").Value = "  +2.64%  "
# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "69.65"
$ws.Range("E18").Value = "  +0.80%  "
# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "253.70"
$ws.Range("E19").Value = "  +0.04%  "
# Row 20
$ws.Range("D20").Value = "0.0₃0800"
$ws.Range("E20").Value = "  +8.00%  "
# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "11.29"
$ws.Range("E21").Value = "  +8.20%  "
# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.997"
$ws.Range("E22").Value = "  -0.38%  "
# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "4.29"
$ws.Range("E23").Value = "  +2.00%  "
# Row 24
$ws.Range("E24").Value = "  +1.36%  "
# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "161.74"
$ws.Range("E25").Value = "  +3.56%  "
# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "16.61"
$ws.Range("E26").Value = "  +0.97%  "
# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "7.21"
$ws.Range("E27").Value = "  +2.50%  "
# Row 28
$ws.Range("E28").Value = "  +0.74%  "
# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.997"
$ws.Range("E29").Value = "  -0.55%  "
# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.0530"
$ws.Range("E30").Value = "  +2.56%  "
# Row 31
$ws.Range("E31").Value = "  +1.12%  "
# Row 32
$ws.Range("E32").Value = "  +0.13%  "
# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "498.64"
$ws.Range("E33").Value = "  +865.12%  "
# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.64"
$ws.Range("E34").Value = "  +2.15%  "
# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.93"
$ws.Range("E35").Value = "  +5.37%  "
# Row 36
$ws.Range("D36").Value = "1.437.25"
$ws.Range("E36").Value = "  -0.91%  "
# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.655"
$ws.Range("E37").Value = "  +3.94%  "
# Row 38
$ws.Range("E38").Value = "  +0.93%  "
# Row 39
$ws.Range("E39").Value = "  +2.50%  "
# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.976"
$ws.Range("E40").Value = "  +9.07%  "
# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "82.78"
$ws.Range("E41").Value = "  -0.39%  "
# Row 42
$ws.Range("E42").Value = "  -2.29%  "
# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.36"
$ws.Range("E43").Value = "  +0.58%  "
# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.17"
$ws.Range("E44").Value = "  +4.90%  "
# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "6.10"
$ws.Range("E45").Value = "  +5.08%  "
# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "12.49"
$ws.Range("E46").Value = "  +5.10%  "
# Row 47
$ws.Range("E47").Value = "  -0.66%  "
# Row 48
$ws.Range("B48").Value = "RocketPoolETH"
$ws.Range("C48").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("D48").Value = "1.975.22"
$ws.Range("E48").Value = "  +1.74%  "
# Row 49
$ws.Range("B49").Value = "Kaspa"
$ws.Range("C49").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0495"
$ws.Range("E49").Value = "  -2.78%  "
# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "106.38"
$ws.Range("E50").Value = "  +8.91%  "
# Row 51
$ws.Range("E51").Value = "  -0.17%  "
